$wb = $excel.ActiveWorkbook

# --- "URL" sheet: rename Profile page -> Home page, add About page / saucelabs URL ---
$wsUrl = $wb.Worksheets.Item("URL")

$wsUrl.Range("B1").Value = "Home page"
$wsUrl.Range("C1").Value = "About page"
$wsUrl.Range("C2").Value = "https://saucelabs.com/"
$wsUrl.Range("C2").Style = "Hyperlink"

# Make "URL" the active sheet/tab, and select cell A3 on it
$wsUrl.Activate() | Out-Null
$wsUrl.Range("A3").Select() | Out-Null
